$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.78%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.39%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08066"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.92%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.887"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.94%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.837"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.77%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9311"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.64%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1284"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-12.39%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1902"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.48%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09240"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.01%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03511"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.29%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09896"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.50%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001434"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006381"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "10.52%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.660"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.50%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.149"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.37%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.164"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.76%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3450"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.60%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.202"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2536"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.42%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004706"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.02%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.10%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-29.60%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05156"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007559"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.67%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01017"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-7.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1370"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.85%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002164"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01077"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.57%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006350"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.96"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.45%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001662"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.45%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
